# Scheduled market-data refresh: update currentAveragePrice* / Leve* columns (H:N)
# across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets with latest pulled values.
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")

# Row 28
$ws.Cells.Item(28, 8).Value = 400.42856  # H28 (currentAveragePrice)
$ws.Cells.Item(28, 9).Value = 400.42856  # I28 (currentAveragePriceNQ)
$ws.Cells.Item(28, 10).Value = 0  # J28 (currentAveragePriceHQ)
$ws.Cells.Item(28, 11).Value = 400.42856  # K28 (LevePriceNQ)
$ws.Cells.Item(28, 12).Value = 0  # L28 (LevePriceHQ)
$ws.Cells.Item(28, 13).Value = 84.57144  # M28 (LeveProfitNQ)
$ws.Cells.Item(28, 14).Value = ""  # N28 (LeveProfitHQ) cleared

# Row 33
$ws.Cells.Item(33, 8).Value = 421.15384  # H33 (currentAveragePrice)
$ws.Cells.Item(33, 9).Value = 408.7  # I33 (currentAveragePriceNQ)
$ws.Cells.Item(33, 11).Value = 408.7  # K33 (LevePriceNQ)
$ws.Cells.Item(33, 13).Value = -179.7  # M33 (LeveProfitNQ)

# Row 69
$ws.Cells.Item(69, 8).Value = 17557  # H69 (currentAveragePrice)
$ws.Cells.Item(69, 9).Value = 3000  # I69 (currentAveragePriceNQ)
$ws.Cells.Item(69, 10).Value = 19983.166  # J69 (currentAveragePriceHQ)
$ws.Cells.Item(69, 11).Value = 9000  # K69 (LevePriceNQ)
$ws.Cells.Item(69, 12).Value = 59949.49800000001  # L69 (LevePriceHQ)
$ws.Cells.Item(69, 13).Value = -8126  # M69 (LeveProfitNQ)
$ws.Cells.Item(69, 14).Value = -61697.49800000001  # N69 (LeveProfitHQ)

# Row 72
$ws.Cells.Item(72, 8).Value = 17557  # H72 (currentAveragePrice)
$ws.Cells.Item(72, 9).Value = 3000  # I72 (currentAveragePriceNQ)
$ws.Cells.Item(72, 10).Value = 19983.166  # J72 (currentAveragePriceHQ)
$ws.Cells.Item(72, 11).Value = 27000  # K72 (LevePriceNQ)
$ws.Cells.Item(72, 12).Value = 179848.494  # L72 (LevePriceHQ)
$ws.Cells.Item(72, 13).Value = -22632  # M72 (LeveProfitNQ)
$ws.Cells.Item(72, 14).Value = -188584.494  # N72 (LeveProfitHQ)

# Row 96
$ws.Cells.Item(96, 8).Value = 1735.5  # H96 (currentAveragePrice)
$ws.Cells.Item(96, 9).Value = 2405.4  # I96 (currentAveragePriceNQ)
$ws.Cells.Item(96, 10).Value = 1257  # J96 (currentAveragePriceHQ)
$ws.Cells.Item(96, 11).Value = 7216.200000000001  # K96 (LevePriceNQ)
$ws.Cells.Item(96, 12).Value = 3771  # L96 (LevePriceHQ)
$ws.Cells.Item(96, 13).Value = -5843.200000000001  # M96 (LeveProfitNQ)
$ws.Cells.Item(96, 14).Value = -6517  # N96 (LeveProfitHQ)

# Row 106
$ws.Cells.Item(106, 8).Value = 11862  # H106 (currentAveragePrice)
$ws.Cells.Item(106, 9).Value = 9270.857  # I106 (currentAveragePriceNQ)
$ws.Cells.Item(106, 10).Value = 30000  # J106 (currentAveragePriceHQ)
$ws.Cells.Item(106, 11).Value = 9270.857  # K106 (LevePriceNQ)
$ws.Cells.Item(106, 12).Value = 30000  # L106 (LevePriceHQ)
$ws.Cells.Item(106, 13).Value = -8639.857  # M106 (LeveProfitNQ)
$ws.Cells.Item(106, 14).Value = -31262  # N106 (LeveProfitHQ)

# Row 137
$ws.Cells.Item(137, 8).Value = 1996  # H137 (currentAveragePrice)
$ws.Cells.Item(137, 10).Value = 1621  # J137 (currentAveragePriceHQ)
$ws.Cells.Item(137, 12).Value = 4863  # L137 (LevePriceHQ)
$ws.Cells.Item(137, 14).Value = -9963  # N137 (LeveProfitHQ)

# Row 138
$ws.Cells.Item(138, 8).Value = 3797.75  # H138 (currentAveragePrice)
$ws.Cells.Item(138, 9).Value = 2998  # I138 (currentAveragePriceNQ)
$ws.Cells.Item(138, 11).Value = 8994  # K138 (LevePriceNQ)
$ws.Cells.Item(138, 13).Value = -3854  # M138 (LeveProfitNQ)

# Row 141
$ws.Cells.Item(141, 8).Value = 8001.8237  # H141 (currentAveragePrice)
$ws.Cells.Item(141, 9).Value = 7573.857  # I141 (currentAveragePriceNQ)
$ws.Cells.Item(141, 11).Value = 22721.571  # K141 (LevePriceNQ)
$ws.Cells.Item(141, 13).Value = -17541.571  # M141 (LeveProfitNQ)


# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")

# Row 5
$ws.Cells.Item(5, 8).Value = 2020  # H5 (currentAveragePrice)
$ws.Cells.Item(5, 10).Value = 524.5  # J5 (currentAveragePriceHQ)
$ws.Cells.Item(5, 12).Value = 524.5  # L5 (LevePriceHQ)
$ws.Cells.Item(5, 14).Value = -748.5  # N5 (LeveProfitHQ)

# Row 32
$ws.Cells.Item(32, 8).Value = 1127.0968  # H32 (currentAveragePrice)
$ws.Cells.Item(32, 9).Value = 1032.4138  # I32 (currentAveragePriceNQ)
$ws.Cells.Item(32, 11).Value = 1032.4138  # K32 (LevePriceNQ)
$ws.Cells.Item(32, 13).Value = -745.4138  # M32 (LeveProfitNQ)

# Row 37
$ws.Cells.Item(37, 8).Value = 4200  # H37 (currentAveragePrice)
$ws.Cells.Item(37, 9).Value = 4200  # I37 (currentAveragePriceNQ)
$ws.Cells.Item(37, 11).Value = 4200  # K37 (LevePriceNQ)
$ws.Cells.Item(37, 13).Value = -3927  # M37 (LeveProfitNQ)

# Row 45
$ws.Cells.Item(45, 8).Value = 2672  # H45 (currentAveragePrice)
$ws.Cells.Item(45, 9).Value = 2515  # I45 (currentAveragePriceNQ)
$ws.Cells.Item(45, 10).Value = 3300  # J45 (currentAveragePriceHQ)
$ws.Cells.Item(45, 11).Value = 2515  # K45 (LevePriceNQ)
$ws.Cells.Item(45, 12).Value = 3300  # L45 (LevePriceHQ)
$ws.Cells.Item(45, 13).Value = -2138  # M45 (LeveProfitNQ)
$ws.Cells.Item(45, 14).Value = -4054  # N45 (LeveProfitHQ)

# Row 102
$ws.Cells.Item(102, 8).Value = 35715984  # H102 (currentAveragePrice)
$ws.Cells.Item(102, 9).Value = 35715984  # I102 (currentAveragePriceNQ)
$ws.Cells.Item(102, 11).Value = 35715984  # K102 (LevePriceNQ)
$ws.Cells.Item(102, 13).Value = -35714362  # M102 (LeveProfitNQ)


# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")

# Row 4
$ws.Cells.Item(4, 8).Value = 2020  # H4 (currentAveragePrice)
$ws.Cells.Item(4, 10).Value = 524.5  # J4 (currentAveragePriceHQ)
$ws.Cells.Item(4, 12).Value = 524.5  # L4 (LevePriceHQ)
$ws.Cells.Item(4, 14).Value = -754.5  # N4 (LeveProfitHQ)

# Row 75
$ws.Cells.Item(75, 8).Value = 5000  # H75 (currentAveragePrice)
$ws.Cells.Item(75, 9).Value = 5000  # I75 (currentAveragePriceNQ)
$ws.Cells.Item(75, 10).Value = 0  # J75 (currentAveragePriceHQ)
$ws.Cells.Item(75, 11).Value = 5000  # K75 (LevePriceNQ)
$ws.Cells.Item(75, 12).Value = 0  # L75 (LevePriceHQ)
$ws.Cells.Item(75, 13).Value = -4064  # M75 (LeveProfitNQ)
$ws.Cells.Item(75, 14).Value = ""  # N75 (LeveProfitHQ) cleared

# Row 78
$ws.Cells.Item(78, 8).Value = 5000  # H78 (currentAveragePrice)
$ws.Cells.Item(78, 9).Value = 5000  # I78 (currentAveragePriceNQ)
$ws.Cells.Item(78, 10).Value = 0  # J78 (currentAveragePriceHQ)
$ws.Cells.Item(78, 11).Value = 15000  # K78 (LevePriceNQ)
$ws.Cells.Item(78, 12).Value = 0  # L78 (LevePriceHQ)
$ws.Cells.Item(78, 13).Value = -10320  # M78 (LeveProfitNQ)
$ws.Cells.Item(78, 14).Value = ""  # N78 (LeveProfitHQ) cleared

# Row 134
$ws.Cells.Item(134, 8).Value = 2564.25  # H134 (currentAveragePrice)
$ws.Cells.Item(134, 9).Value = 1500  # I134 (currentAveragePriceNQ)
$ws.Cells.Item(134, 11).Value = 4500  # K134 (LevePriceNQ)
$ws.Cells.Item(134, 13).Value = -1965  # M134 (LeveProfitNQ)


# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")

# Row 16
$ws.Cells.Item(16, 8).Value = 1790.75  # H16 (currentAveragePrice)
$ws.Cells.Item(16, 9).Value = 1688  # I16 (currentAveragePriceNQ)
$ws.Cells.Item(16, 11).Value = 1688  # K16 (LevePriceNQ)
$ws.Cells.Item(16, 13).Value = -1401  # M16 (LeveProfitNQ)

# Row 31
$ws.Cells.Item(31, 8).Value = 1875.2  # H31 (currentAveragePrice)
$ws.Cells.Item(31, 9).Value = 1875.2  # I31 (currentAveragePriceNQ)
$ws.Cells.Item(31, 11).Value = 1875.2  # K31 (LevePriceNQ)
$ws.Cells.Item(31, 13).Value = -1580.2  # M31 (LeveProfitNQ)

# Row 34
$ws.Cells.Item(34, 8).Value = 1875.2  # H34 (currentAveragePrice)
$ws.Cells.Item(34, 9).Value = 1875.2  # I34 (currentAveragePriceNQ)
$ws.Cells.Item(34, 11).Value = 1875.2  # K34 (LevePriceNQ)
$ws.Cells.Item(34, 13).Value = -1673.2  # M34 (LeveProfitNQ)

# Row 63
$ws.Cells.Item(63, 8).Value = 100271  # H63 (currentAveragePrice)
$ws.Cells.Item(63, 10).Value = 100271  # J63 (currentAveragePriceHQ)
$ws.Cells.Item(63, 12).Value = 100271  # L63 (LevePriceHQ)
$ws.Cells.Item(63, 14).Value = -101643  # N63 (LeveProfitHQ)

# Row 66
$ws.Cells.Item(66, 8).Value = 100271  # H66 (currentAveragePrice)
$ws.Cells.Item(66, 10).Value = 100271  # J66 (currentAveragePriceHQ)
$ws.Cells.Item(66, 12).Value = 300813  # L66 (LevePriceHQ)
$ws.Cells.Item(66, 14).Value = -307677  # N66 (LeveProfitHQ)

# Row 113
$ws.Cells.Item(113, 8).Value = 1790.75  # H113 (currentAveragePrice)
$ws.Cells.Item(113, 9).Value = 1688  # I113 (currentAveragePriceNQ)
$ws.Cells.Item(113, 11).Value = 1688  # K113 (LevePriceNQ)
$ws.Cells.Item(113, 13).Value = 482  # M113 (LeveProfitNQ)

# Row 125
$ws.Cells.Item(125, 8).Value = 0  # H125 (currentAveragePrice)
$ws.Cells.Item(125, 10).Value = 0  # J125 (currentAveragePriceHQ)
$ws.Cells.Item(125, 12).Value = 0  # L125 (LevePriceHQ)
$ws.Cells.Item(125, 14).Value = ""  # N125 (LeveProfitHQ) cleared


# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")

# Row 2
$ws.Cells.Item(2, 8).Value = 92.27273  # H2 (currentAveragePrice)
$ws.Cells.Item(2, 10).Value = 99.625  # J2 (currentAveragePriceHQ)
$ws.Cells.Item(2, 12).Value = 597.75  # L2 (LevePriceHQ)
$ws.Cells.Item(2, 14).Value = -823.75  # N2 (LeveProfitHQ)

# Row 98
$ws.Cells.Item(98, 8).Value = 613.75  # H98 (currentAveragePrice)
$ws.Cells.Item(98, 10).Value = 718.3333  # J98 (currentAveragePriceHQ)
$ws.Cells.Item(98, 12).Value = 2154.9999  # L98 (LevePriceHQ)
$ws.Cells.Item(98, 14).Value = -5150.9999  # N98 (LeveProfitHQ)

# Row 120
$ws.Cells.Item(120, 8).Value = 0  # H120 (currentAveragePrice)
$ws.Cells.Item(120, 9).Value = 0  # I120 (currentAveragePriceNQ)
$ws.Cells.Item(120, 11).Value = 0  # K120 (LevePriceNQ)
$ws.Cells.Item(120, 13).Value = ""  # M120 (LeveProfitNQ) cleared


# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")

# Row 80
$ws.Cells.Item(80, 8).Value = 1565.1765  # H80 (currentAveragePrice)
$ws.Cells.Item(80, 9).Value = 1418.6666  # I80 (currentAveragePriceNQ)
$ws.Cells.Item(80, 10).Value = 1916.8  # J80 (currentAveragePriceHQ)
$ws.Cells.Item(80, 11).Value = 1418.6666  # K80 (LevePriceNQ)
$ws.Cells.Item(80, 12).Value = 1916.8  # L80 (LevePriceHQ)
$ws.Cells.Item(80, 13).Value = -420.6666  # M80 (LeveProfitNQ)
$ws.Cells.Item(80, 14).Value = -3912.8  # N80 (LeveProfitHQ)

# Row 83
$ws.Cells.Item(83, 8).Value = 1565.1765  # H83 (currentAveragePrice)
$ws.Cells.Item(83, 9).Value = 1418.6666  # I83 (currentAveragePriceNQ)
$ws.Cells.Item(83, 10).Value = 1916.8  # J83 (currentAveragePriceHQ)
$ws.Cells.Item(83, 11).Value = 7093.333000000001  # K83 (LevePriceNQ)
$ws.Cells.Item(83, 12).Value = 9584  # L83 (LevePriceHQ)
$ws.Cells.Item(83, 13).Value = -2101.333000000001  # M83 (LeveProfitNQ)
$ws.Cells.Item(83, 14).Value = -19568  # N83 (LeveProfitHQ)

# Row 102
$ws.Cells.Item(102, 8).Value = 2642.1875  # H102 (currentAveragePrice)
$ws.Cells.Item(102, 9).Value = 2718.3333  # I102 (currentAveragePriceNQ)
$ws.Cells.Item(102, 10).Value = 1500  # J102 (currentAveragePriceHQ)
$ws.Cells.Item(102, 11).Value = 2718.3333  # K102 (LevePriceNQ)
$ws.Cells.Item(102, 12).Value = 1500  # L102 (LevePriceHQ)
$ws.Cells.Item(102, 13).Value = -1096.3333  # M102 (LeveProfitNQ)
$ws.Cells.Item(102, 14).Value = -4744  # N102 (LeveProfitHQ)

# Row 122
$ws.Cells.Item(122, 8).Value = 2636.8572  # H122 (currentAveragePrice)
$ws.Cells.Item(122, 9).Value = 1983.8572  # I122 (currentAveragePriceNQ)
$ws.Cells.Item(122, 11).Value = 5951.571599999999  # K122 (LevePriceNQ)
$ws.Cells.Item(122, 13).Value = -3501.571599999999  # M122 (LeveProfitNQ)

# Row 132
$ws.Cells.Item(132, 8).Value = 5386.125  # H132 (currentAveragePrice)
$ws.Cells.Item(132, 9).Value = 5014.8335  # I132 (currentAveragePriceNQ)
$ws.Cells.Item(132, 11).Value = 15044.5005  # K132 (LevePriceNQ)
$ws.Cells.Item(132, 13).Value = -12514.5005  # M132 (LeveProfitNQ)


# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")

# Row 46
$ws.Cells.Item(46, 8).Value = 3186.75  # H46 (currentAveragePrice)
$ws.Cells.Item(46, 9).Value = 2666.6667  # I46 (currentAveragePriceNQ)
$ws.Cells.Item(46, 10).Value = 3498.8  # J46 (currentAveragePriceHQ)
$ws.Cells.Item(46, 11).Value = 2666.6667  # K46 (LevePriceNQ)
$ws.Cells.Item(46, 12).Value = 3498.8  # L46 (LevePriceHQ)
$ws.Cells.Item(46, 13).Value = -2478.6667  # M46 (LeveProfitNQ)
$ws.Cells.Item(46, 14).Value = -3874.8  # N46 (LeveProfitHQ)

# Row 48
$ws.Cells.Item(48, 8).Value = 33360.668  # H48 (currentAveragePrice)
$ws.Cells.Item(48, 9).Value = 33360.668  # I48 (currentAveragePriceNQ)
$ws.Cells.Item(48, 11).Value = 33360.668  # K48 (LevePriceNQ)
$ws.Cells.Item(48, 13).Value = -32699.668  # M48 (LeveProfitNQ)

# Row 53
$ws.Cells.Item(53, 8).Value = 16950  # H53 (currentAveragePrice)
$ws.Cells.Item(53, 9).Value = 17000  # I53 (currentAveragePriceNQ)
$ws.Cells.Item(53, 10).Value = 16900  # J53 (currentAveragePriceHQ)
$ws.Cells.Item(53, 11).Value = 17000  # K53 (LevePriceNQ)
$ws.Cells.Item(53, 12).Value = 16900  # L53 (LevePriceHQ)
$ws.Cells.Item(53, 13).Value = -16482  # M53 (LeveProfitNQ)
$ws.Cells.Item(53, 14).Value = -17936  # N53 (LeveProfitHQ)

# Row 61
$ws.Cells.Item(61, 8).Value = 1440.1305  # H61 (currentAveragePrice)
$ws.Cells.Item(61, 9).Value = 1291.1666  # I61 (currentAveragePriceNQ)
$ws.Cells.Item(61, 10).Value = 1976.4  # J61 (currentAveragePriceHQ)
$ws.Cells.Item(61, 11).Value = 1291.1666  # K61 (LevePriceNQ)
$ws.Cells.Item(61, 12).Value = 1976.4  # L61 (LevePriceHQ)
$ws.Cells.Item(61, 13).Value = -1089.1666  # M61 (LeveProfitNQ)
$ws.Cells.Item(61, 14).Value = -2380.4  # N61 (LeveProfitHQ)

# Row 93
$ws.Cells.Item(93, 8).Value = 2175  # H93 (currentAveragePrice)
$ws.Cells.Item(93, 9).Value = 1066.6666  # I93 (currentAveragePriceNQ)
$ws.Cells.Item(93, 11).Value = 1066.6666  # K93 (LevePriceNQ)
$ws.Cells.Item(93, 13).Value = 181.3334  # M93 (LeveProfitNQ)

# Row 94
$ws.Cells.Item(94, 8).Value = 29999  # H94 (currentAveragePrice)
$ws.Cells.Item(94, 10).Value = 29999  # J94 (currentAveragePriceHQ)
$ws.Cells.Item(94, 12).Value = 29999  # L94 (LevePriceHQ)
$ws.Cells.Item(94, 14).Value = -31351  # N94 (LeveProfitHQ)

# Row 113
$ws.Cells.Item(113, 8).Value = 1440.1305  # H113 (currentAveragePrice)
$ws.Cells.Item(113, 9).Value = 1291.1666  # I113 (currentAveragePriceNQ)
$ws.Cells.Item(113, 10).Value = 1976.4  # J113 (currentAveragePriceHQ)
$ws.Cells.Item(113, 11).Value = 1291.1666  # K113 (LevePriceNQ)
$ws.Cells.Item(113, 12).Value = 1976.4  # L113 (LevePriceHQ)
$ws.Cells.Item(113, 13).Value = 878.8334  # M113 (LeveProfitNQ)
$ws.Cells.Item(113, 14).Value = -6316.4  # N113 (LeveProfitHQ)

# Row 136
$ws.Cells.Item(136, 8).Value = 7484.25  # H136 (currentAveragePrice)
$ws.Cells.Item(136, 9).Value = 9314.333000000001  # I136 (currentAveragePriceNQ)
$ws.Cells.Item(136, 10).Value = 1994  # J136 (currentAveragePriceHQ)
$ws.Cells.Item(136, 11).Value = 27942.999  # K136 (LevePriceNQ)
$ws.Cells.Item(136, 12).Value = 5982  # L136 (LevePriceHQ)
$ws.Cells.Item(136, 13).Value = -25392.999  # M136 (LeveProfitNQ)
$ws.Cells.Item(136, 14).Value = -11082  # N136 (LeveProfitHQ)


# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")

# Row 82
$ws.Cells.Item(82, 8).Value = 20000  # H82 (currentAveragePrice)
$ws.Cells.Item(82, 10).Value = 20000  # J82 (currentAveragePriceHQ)
$ws.Cells.Item(82, 12).Value = 20000  # L82 (LevePriceHQ)
$ws.Cells.Item(82, 14).Value = -20766  # N82 (LeveProfitHQ)

# Row 85
$ws.Cells.Item(85, 8).Value = 20000  # H85 (currentAveragePrice)
$ws.Cells.Item(85, 10).Value = 20000  # J85 (currentAveragePriceHQ)
$ws.Cells.Item(85, 12).Value = 20000  # L85 (LevePriceHQ)
$ws.Cells.Item(85, 14).Value = -22652  # N85 (LeveProfitHQ)

# Row 113
$ws.Cells.Item(113, 8).Value = 1024.5555  # H113 (currentAveragePrice)
$ws.Cells.Item(113, 9).Value = 1065.25  # I113 (currentAveragePriceNQ)
$ws.Cells.Item(113, 10).Value = 699  # J113 (currentAveragePriceHQ)
$ws.Cells.Item(113, 11).Value = 3195.75  # K113 (LevePriceNQ)
$ws.Cells.Item(113, 12).Value = 2097  # L113 (LevePriceHQ)
$ws.Cells.Item(113, 13).Value = -1025.75  # M113 (LeveProfitNQ)
$ws.Cells.Item(113, 14).Value = -6437  # N113 (LeveProfitHQ)
